# Branch - Brand - SKU wise Stock Aging Status: Summary
# Re-order the Item Name / UOM pairs within brand groups so the shared
# string table (and therefore the cell values) reflect the updated order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dinafex group (rows 3-5): 60mg, 180mg, 120mg -> 120mg, 60mg, 180mg
$ws.Range("D3").Value = "Dinafex 120mg Tablet"
$ws.Range("D4").Value = "Dinafex 60mg Tablet"
$ws.Range("D5").Value = "Dinafex 180mg Tablet"

# Etorix group (rows 7-9): 120mg, 60mg-40's, 90mg -> 60mg-40's, 90mg, 120mg
$ws.Range("D7").Value = "Etorix 60mg Tablet - 40's"
$ws.Range("E7").Value = "40's"
$ws.Range("D8").Value = "Etorix 90mg Tablet"
$ws.Range("E8").Value = "30's"
$ws.Range("D9").Value = "Etorix 120mg Tablet"
$ws.Range("E9").Value = "20's"

# Ketonic tablet/injection pair (rows 15-16): 10mg Tablet, IM/IV 4's -> IM/IV 4's, 10mg Tablet
$ws.Range("D15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E15").Value = "4's"
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

# Kynol group (rows 18-19): D 25mg, TR 100mg -> TR 100mg, D 25mg
$ws.Range("D18").Value = "Kynol TR 100mg Capsule"
$ws.Range("E18").Value = "50 's"
$ws.Range("D19").Value = "Kynol D 25mg Tablet"
$ws.Range("E19").Value = "60 's"

# Zithrox group (rows 24-27): 250mg-6's, 15ml Susp, 500mg, 30ml Dry -> 500mg, 30ml Dry, 250mg-6's, 15ml Susp
$ws.Range("D24").Value = "Zithrox 500mg Tablet"
$ws.Range("E24").Value = "6 's"
$ws.Range("D25").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E25").Value = "30ml"
$ws.Range("D26").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E26").Value = "6's"
$ws.Range("D27").Value = "Zithrox 15ml Suspension"
$ws.Range("E27").Value = "15 ml"
